$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Common Word"
$ws.Range("B1").Value = "Total Frequency"
$ws.Range("C1").Value = "Websites"

# Data rows
$ws.Range("A2").Value = "A3"
$ws.Range("B2").Value = 71
$ws.Range("C2").Value = "http://www.yyy.at/ (54), https://www.xxx.at/ (6), https://www.zzz.at/ (11)"

$ws.Range("A3").Value = "B3"
$ws.Range("B3").Value = 68
$ws.Range("C3").Value = "http://www.yyy.at/ (52), https://www.xxx.at/ (6), https://www.zzz.at/ (10)"

$ws.Range("A4").Value = "C2"
$ws.Range("B4").Value = 55
$ws.Range("C4").Value = "http://www.yyy.at/ (50), https://www.xxx.at/ (5)"

$ws.Range("A5").Value = "D2"
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = "https://www.xxx.at/ (5), https://www.zzz.at/ (8)"

# Header formatting: bold font, thin box border, centered horizontally, top vertically
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108  # xlCenter
$header.VerticalAlignment = -4160    # xlTop

$header.Borders.Item(7).LineStyle = 1   # xlEdgeLeft, xlContinuous
$header.Borders.Item(7).Weight = 2      # xlThin
$header.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$header.Borders.Item(8).Weight = 2
$header.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$header.Borders.Item(9).Weight = 2
$header.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$header.Borders.Item(10).Weight = 2
